# LOB1012.xlsx edit
# 1) Replace the "Objetivos" description (B10/C10) with the docente text that used
#    to live in row 13.
# 2) Delete row 13 (the old "Docentes responsáveis" value row, now merged one row
#    up with "Programa resumido:"), shifting everything below up by one row.
# 3) Re-populate the cells whose text content changed as part of the shuffle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Objetivos value becomes the docente identification string.
$ws.Range("B10").Value = "4894221 - Mariana Pereira de Melo"
$ws.Range("C10").Value = "4894221 - Mariana Pereira de Melo"

# 2) Remove row 13 entirely; rows 14-24 shift up to 13-23.
$ws.Rows(13).Delete()

# 3) Fix up the content of the rows affected by the shift.

# Row 13 (was "Programa resumido:" / long description) -> now "Semestral"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (was "Programa:" / long syllabus) -> now the activation date string,
# reusing the existing "01/01/2018" text (copy as value to avoid Excel
# re-interpreting the literal string as a date).
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4163) | Out-Null
$ws.Range("C8").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# Row 18 (was "Avaliação:", blank B/C) -> now "Método:" paired with the docente string.
$ws.Range("B18").Value = "4894221 - Mariana Pereira de Melo"
$ws.Range("C18").Value = "4894221 - Mariana Pereira de Melo"

# Row 19 (was "Método:" content) -> now "Critério:" paired with the NF formula text.
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# Row 20 (was "Critério:" content) -> now "Norma de recuperação:" paired with "NF≥ 5,0."
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

# Row 21 (was "Norma de recuperação:" content) -> now "Bibliografia:" paired with the
# recovery-exam criterion text.
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
